# Adds season-record columns (Wins, Losses, Ties) to the CHC_2017 roster
# sheet: new header cells AD1:AF1 and a constant 92-70-0 record copied down
# every existing data row (rows 2-49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Copy the formatting of the existing last header cell (AC1, style "1":
# bold, bordered, centered) into the three new header cells first, then
# overwrite their text - this keeps them on the same cell style as the
# rest of row 1 instead of minting a brand-new style.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows ----------------------------------------------------------
# Every player row gets the team's 2017 season record: 92 wins, 70
# losses, 0 ties.
for ($r = 2; $r -le 49; $r++) {
    $ws.Range("AD$r").Value = 92
    $ws.Range("AE$r").Value = 70
    $ws.Range("AF$r").Value = 0
}
